{"js": "const body = context.document.body;\n\n// Map of exact old text -> new text. Each old value is unique in the\n// document, so we issue ALL searches first (against the original,\n// unmodified content) and only perform the text replacements afterwards.\n// This avoids a replacement's new text accidentally being re-matched by a\n// later search (e.g. '48\u00f79=' -> '76\u00f76=' followed by the original '76\u00f76=' -> '97\u00f79=').\nconst replacements = [\n  { find: \"2024-07-12 Friday\", replace: \"2024-07-13 Saturday\" },\n  { find: \"67\u00f72=\", replace: \"51\u00f76=\" },\n  { find: \"68\u00f77=\", replace: \"17\u00f73=\" },\n  { find: \"88\u00f76=\", replace: \"28\u00f78=\" },\n  { find: \"48\u00f79=\", replace: \"76\u00f76=\" },\n  { find: \"76\u00f78=\", replace: \"59\u00f79=\" },\n  { find: \"76\u00f76=\", replace: \"97\u00f79=\" },\n  { find: \"93\u00f76=\", replace: \"19\u00f73=\" },\n  { find: \"31\u00f72=\", replace: \"81\u00f77=\" },\n  { find: \"63\u00f75=\", replace: \"78\u00f78=\" },\n  { find: \"47\u00f78=\", replace: \"78\u00f73=\" },\n  { find: \"67\u00f77=\", replace: \"52\u00f78=\" },\n  { find: \"78\u00f79=\", replace: \"81\u00f72=\" },\n  { find: \"79\u00f76=\", replace: \"32\u00f76=\" },\n  { find: \"84\u00f78=\", replace: \"33\u00f74=\" },\n  { find: \"66\u00f74=\", replace: \"42\u00f76=\" },\n  { find: \"77\u00f77=\", replace: \"14\u00f76=\" },\n  { find: \"31\u00f76=\", replace: \"66\u00f74=\" },\n  { find: \"30\u00f79=\", replace: \"83\u00f74=\" },\n  { find: \"19\u00f74=\", replace: \"96\u00f79=\" },\n  { find: \"98\u00f75=\", replace: \"43\u00f78=\" },\n  { find: \"69\u00f76=\", replace: \"74\u00f75=\" },\n  { find: \"99\u00f78=\", replace: \"80\u00f75=\" },\n  { find: \"85\u00f78=\", replace: \"41\u00f75=\" },\n  { find: \"63\u00f72=\", replace: \"88\u00f73=\" },\n  { find: \"98\u00f77=\", replace: \"68\u00f74=\" },\n];\n\nconst searchResults = replacements.map(r =>\n  body.search(r.find, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach(r => r.load('items'));\nawait context.sync();\n\nsearchResults.forEach((r, i) => {\n  if (r.items.length === 0) {\n    throw new Error('Text not found: ' + replacements[i].find);\n  }\n  r.items[0].insertText(replacements[i].replace, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph of the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-07-13 Saturday\"\n\n# Update each division-problem cell by explicit (row, column) coordinates\n# in the table, which sidesteps any ambiguity from duplicate text values\n# (e.g. '76\u00f76=' appears twice before the edit, at two different cells,\n# and changes to two different results) and avoids re-matching text that\n# a previous replacement just inserted.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n  1 = @(\"51\u00f76=\", \"17\u00f73=\", \"28\u00f78=\", \"76\u00f76=\", \"59\u00f79=\")\n  5 = @(\"97\u00f79=\", \"19\u00f73=\", \"81\u00f77=\", \"78\u00f78=\", \"78\u00f73=\")\n  9 = @(\"52\u00f78=\", \"81\u00f72=\", \"32\u00f76=\", \"33\u00f74=\", \"42\u00f76=\")\n  13 = @(\"14\u00f76=\", \"66\u00f74=\", \"83\u00f74=\", \"96\u00f79=\", \"43\u00f78=\")\n  17 = @(\"74\u00f75=\", \"80\u00f75=\", \"41\u00f75=\", \"88\u00f73=\", \"68\u00f74=\")\n}\n\nforeach ($row in $newValues.Keys) {\n  $rowValues = $newValues[$row]\n  for ($col = 1; $col -le $rowValues.Count; $col++) {\n    $t.Cell($row, $col).Range.Text = $rowValues[$col - 1]\n  }\n}\n"}
